$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- "Ready for handoff" -> "In Translation" (the shared status string, used on
#     every sheet: Overview!E/F2:E/F3 and the per-locale Status column C2:C3) ---
$newStatus = "In Translation"

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# --- Narrow the "Status" columns: Overview columns E & F, and column C on the
#     zh-cn / de-de sheets (was width 17.2159881591797, now 13.4101848602295) ---
# A ColumnWidth (character units) of 12.5 renders, once Excel applies its
# standard cell padding, to the narrower stored column width used in the sheet.
$newColumnWidth = 12.5

$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth   # column E
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth   # column F

$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth       # column C

$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth       # column C
